# Update "Valori stimati MLP" (column G, rows 2-19) on sheets Caso1, Caso2, Caso3
# with newly uploaded MLP model results.

$wb = $excel.ActiveWorkbook

$caso1Values = @(
    0.9981211423873901,
    0.9980431199073792,
    0.9966030120849609,
    1.00377357006073,
    1.005973339080811,
    0.9934861660003662,
    0.9885492920875549,
    0.9974955916404724,
    0.995840311050415,
    0.9958052039146423,
    0.9937840104103088,
    0.9933983683586121,
    0.9999719262123108,
    1.000558137893677,
    1.004407525062561,
    0.9897051453590393,
    0.9950308799743652,
    0.9936397075653076
)

$caso2Values = @(
    1.000048398971558,
    0.9999738335609436,
    0.9983891248703003,
    1.005804777145386,
    1.007598400115967,
    0.9950535297393799,
    0.990659773349762,
    0.9989772439002991,
    0.9976685047149658,
    0.9977731108665466,
    0.9952053427696228,
    0.995067834854126,
    1.001819491386414,
    1.002386450767517,
    1.006364822387695,
    0.9914928078651428,
    0.9968952536582947,
    0.9953093528747559
)

$caso3Values = @(
    0.999210000038147,
    0.998903214931488,
    0.9976637959480286,
    1.005215048789978,
    1.006742119789124,
    0.9936344623565674,
    0.9898490309715271,
    0.9983139634132385,
    0.9968945980072021,
    0.9968197345733643,
    0.9941394925117493,
    0.9945316314697266,
    1.000686764717102,
    1.001368045806885,
    1.005790591239929,
    0.9904484748840332,
    0.9961087107658386,
    0.9946746826171875
)

$ws1 = $wb.Worksheets.Item("Caso1")
for ($i = 0; $i -lt $caso1Values.Length; $i++) {
    $ws1.Cells.Item(2 + $i, 7).Value = $caso1Values[$i]
}

$ws2 = $wb.Worksheets.Item("Caso2")
for ($i = 0; $i -lt $caso2Values.Length; $i++) {
    $ws2.Cells.Item(2 + $i, 7).Value = $caso2Values[$i]
}

$ws3 = $wb.Worksheets.Item("Caso3")
for ($i = 0; $i -lt $caso3Values.Length; $i++) {
    $ws3.Cells.Item(2 + $i, 7).Value = $caso3Values[$i]
}
